$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Cells.Item(2, 4) "22.340.08"
Set-TextValue $ws.Cells.Item(2, 5) "  -4.87%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.568.86"
Set-TextValue $ws.Cells.Item(3, 5) "  -4.92%  "
Set-TextValue $ws.Cells.Item(4, 5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(5, 5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(6, 4) "289.54"
Set-TextValue $ws.Cells.Item(6, 5) "  -3.52%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.3755"
Set-TextValue $ws.Cells.Item(7, 5) "  -0.66%  "
Set-TextValue $ws.Cells.Item(8, 4) "49.43"
Set-TextValue $ws.Cells.Item(8, 5) "  -2.76%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.3433"
Set-TextValue $ws.Cells.Item(9, 5) "  -1.76%  "
Set-TextValue $ws.Cells.Item(10, 4) "1.170"
Set-TextValue $ws.Cells.Item(10, 5) "  -4.69%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07670"
Set-TextValue $ws.Cells.Item(11, 5) "  -4.80%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.9999"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.08%  "
Set-TextValue $ws.Cells.Item(13, 4) "21.42"
Set-TextValue $ws.Cells.Item(13, 5) "  -2.93%  "
Set-TextValue $ws.Cells.Item(14, 4) "6.032"
Set-TextValue $ws.Cells.Item(14, 5) "  -4.52%  "
Set-TextValue $ws.Cells.Item(15, 4) "6.959"
Set-TextValue $ws.Cells.Item(15, 5) "  -4.23%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.00001138"
Set-TextValue $ws.Cells.Item(16, 5) "  -5.87%  "
Set-TextValue $ws.Cells.Item(17, 4) "1.569.94"
Set-TextValue $ws.Cells.Item(17, 5) "  -5.16%  "
Set-TextValue $ws.Cells.Item(18, 4) "90.24"
Set-TextValue $ws.Cells.Item(18, 5) "  -5.23%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.06720"
Set-TextValue $ws.Cells.Item(19, 5) "  -3.69%  "
Set-TextValue $ws.Cells.Item(20, 5) "  -0.05%  "
Set-TextValue $ws.Cells.Item(21, 4) "6.258"
Set-TextValue $ws.Cells.Item(21, 5) "  -5.80%  "
Set-TextValue $ws.Cells.Item(22, 4) "16.64"
Set-TextValue $ws.Cells.Item(22, 5) "  -4.87%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.5272"
Set-TextValue $ws.Cells.Item(23, 5) "  -8.65%  "
Set-TextValue $ws.Cells.Item(24, 4) "11.97"
Set-TextValue $ws.Cells.Item(24, 5) "  -4.11%  "
Set-TextValue $ws.Cells.Item(25, 4) "22.332.13"
Set-TextValue $ws.Cells.Item(25, 5) "  -4.90%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.391"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.34%  "
Set-TextValue $ws.Cells.Item(27, 4) "2.799"
Set-TextValue $ws.Cells.Item(27, 5) "  -7.46%  "
Set-TextValue $ws.Cells.Item(28, 4) "20.16"
Set-TextValue $ws.Cells.Item(28, 5) "  -4.57%  "
Set-TextValue $ws.Cells.Item(29, 4) "145.41"
Set-TextValue $ws.Cells.Item(29, 5) "  -4.10%  "
Set-TextValue $ws.Cells.Item(30, 4) "4.981"
Set-TextValue $ws.Cells.Item(30, 5) "  -3.80%  "
Set-TextValue $ws.Cells.Item(31, 4) "125.72"
Set-TextValue $ws.Cells.Item(31, 5) "  -4.60%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.734.55"
Set-TextValue $ws.Cells.Item(32, 5) "  -5.63%  "
Set-TextValue $ws.Cells.Item(33, 4) "1.027"
Set-TextValue $ws.Cells.Item(33, 5) "  +3.67%  "
Set-TextValue $ws.Cells.Item(34, 4) "6.238"
Set-TextValue $ws.Cells.Item(34, 5) "  -9.64%  "
Set-TextValue $ws.Cells.Item(35, 4) "2.009"
Set-TextValue $ws.Cells.Item(35, 5) "  -6.34%  "
Set-TextValue $ws.Cells.Item(36, 4) "10.12"
Set-TextValue $ws.Cells.Item(36, 5) "  -8.97%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.08545"
Set-TextValue $ws.Cells.Item(37, 5) "  -2.87%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.02542"
Set-TextValue $ws.Cells.Item(38, 5) "  -7.09%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.2330"
Set-TextValue $ws.Cells.Item(39, 5) "  -3.85%  "
Set-TextValue $ws.Cells.Item(40, 4) "5.552"
Set-TextValue $ws.Cells.Item(40, 5) "  -6.47%  "
Set-TextValue $ws.Cells.Item(41, 4) "1.325"
Set-TextValue $ws.Cells.Item(41, 5) "  +2.00%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.06415"
Set-TextValue $ws.Cells.Item(42, 5) "  -6.53%  "
Set-TextValue $ws.Cells.Item(43, 4) "11.78"
Set-TextValue $ws.Cells.Item(43, 5) "  -8.40%  "
Set-TextValue $ws.Cells.Item(44, 4) "0.6433"
Set-TextValue $ws.Cells.Item(44, 5) "  -6.99%  "
Set-TextValue $ws.Cells.Item(45, 4) "14.27"
Set-TextValue $ws.Cells.Item(45, 5) "  -8.79%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.9995"
Set-TextValue $ws.Cells.Item(46, 5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.6012"
Set-TextValue $ws.Cells.Item(47, 5) "  -6.16%  "
Set-TextValue $ws.Cells.Item(48, 4) "3.755"
Set-TextValue $ws.Cells.Item(48, 5) "  -4.20%  "
Set-TextValue $ws.Cells.Item(49, 4) "2.101"
Set-TextValue $ws.Cells.Item(49, 5) "  -6.70%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.273"
Set-TextValue $ws.Cells.Item(50, 5) "  +2.65%  "
Set-TextValue $ws.Cells.Item(51, 4) "124.43"
Set-TextValue $ws.Cells.Item(51, 5) "  -2.01%  "
